# Fill in the previously-blank row 18 of the time sheet with a new entry:
# Simon worked 1 hour on "Editting files" on 43362 (19/09/2018), matching
# the style/formatting already used by the other "Simon" rows (e.g. row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from row 13 (an existing Simon entry) onto row 18 so
# the new row picks up the same fill/border/number-format as its peers.
$ws.Range("A13:D13").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Now populate the actual values for the new entry.
$ws.Range("A18").Value = "Simon"
$ws.Range("B18").Value = "Editting files"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 43362

# Move the active selection to B18, matching the saved workbook state.
$ws.Range("B18").Select() | Out-Null
